$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the cryptos list refresh
$ws.Range("D2").Value = '27.575.04'
$ws.Range("E2").Value = '  -1.26%  '
$ws.Range("D3").Value = '1.847.19'
$ws.Range("E3").Value = '  -2.11%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -1.39%  '
$ws.Range("D5").Value = '332.27'
$ws.Range("E5").Value = '  -1.05%  '
$ws.Range("E6").Value = '  -1.29%  '
$ws.Range("D7").Value = '0.4625'
$ws.Range("E7").Value = '  -1.93%  '
$ws.Range("D8").Value = '0.3850'
$ws.Range("E8").Value = '  -2.21%  '
$ws.Range("D9").Value = '45.92'
$ws.Range("E9").Value = '  -2.19%  '
$ws.Range("D10").Value = '0.07909'
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("D11").Value = '0.9926'
$ws.Range("E11").Value = '  -2.30%  '
$ws.Range("D12").Value = '21.44'
$ws.Range("E12").Value = '  -1.39%  '
$ws.Range("D13").Value = '1.861.65'
$ws.Range("E13").Value = '  -1.79%  '
$ws.Range("D14").Value = '5.906'
$ws.Range("E14").Value = '  -1.43%  '
$ws.Range("D15").Value = '7.088'
$ws.Range("E15").Value = '  -1.20%  '
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").Value = '  -1.33%  '
$ws.Range("D17").Value = '88.64'
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("D18").Value = '0.06653'
$ws.Range("E18").Value = '  -2.01%  '
$ws.Range("E19").Value = '  -1.61%  '
$ws.Range("D20").Value = '17.04'
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  -1.24%  '
$ws.Range("D22").Value = '27.579.84'
$ws.Range("E22").Value = '  -1.31%  '
$ws.Range("D23").Value = '5.373'
$ws.Range("E23").Value = '  -2.23%  '
$ws.Range("D24").Value = '10.90'
$ws.Range("D25").Value = '2.305'
$ws.Range("E25").Value = '  -2.22%  '
$ws.Range("D26").Value = '157.81'
$ws.Range("E26").Value = '  -0.96%  '
$ws.Range("D27").Value = '19.47'
$ws.Range("E27").Value = '  -2.64%  '
$ws.Range("D28").Value = '2.087'
$ws.Range("E28").Value = '  -0.62%  '
$ws.Range("D29").Value = '5.389'
$ws.Range("E29").Value = '  -1.77%  '
$ws.Range("D30").Value = '119.59'
$ws.Range("E30").Value = '  -1.60%  '
$ws.Range("D31").Value = '0.9718'
$ws.Range("E31").Value = '  +0.80%  '
$ws.Range("D32").Value = '0.09377'
$ws.Range("E32").Value = '  -2.00%  '
$ws.Range("D33").Value = '3.578'
$ws.Range("E33").Value = '  -1.86%  '
$ws.Range("D34").Value = '5.272'
$ws.Range("E34").Value = '  -1.44%  '
$ws.Range("D35").Value = '1.344'
$ws.Range("E35").Value = '  -1.07%  '
$ws.Range("E36").Value = '  -1.85%  '
$ws.Range("D37").Value = '0.02217'
$ws.Range("E37").Value = '  -1.17%  '
$ws.Range("D38").Value = '8.272'
$ws.Range("E38").Value = '  +0.69%  '
$ws.Range("E39").Value = '  -2.89%  '
$ws.Range("D40").Value = '0.5880'
$ws.Range("E40").Value = '  -1.32%  '
$ws.Range("E41").Value = '  -2.46%  '
$ws.Range("E42").Value = '  -0.76%  '
$ws.Range("E43").Value = '  -2.29%  '
$ws.Range("D44").Value = '0.5572'
$ws.Range("E44").Value = '  -1.95%  '
$ws.Range("D45").Value = '12.16'
$ws.Range("E45").Value = '  -0.60%  '
$ws.Range("E46").Value = '  -2.55%  '
$ws.Range("D47").Value = '0.06681'
$ws.Range("E47").Value = '  -2.75%  '
$ws.Range("D48").Value = '110.48'
$ws.Range("E48").Value = '  -2.43%  '
$ws.Range("D49").Value = '1.051'
$ws.Range("E49").Value = '  -1.71%  '
$ws.Range("E50").Value = '  -1.24%  '
$ws.Range("D51").Value = '70.06'
$ws.Range("E51").Value = '  -1.59%  '

Write-Output "done"
